$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.204.17"
$ws.Range("E2").Value = "  +0.87%  "

# Row 3
$ws.Range("D3").Value = "1.852.44"
$ws.Range("E3").Value = "  +1.47%  "

# Row 4
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").Value = "'313.52"
$ws.Range("E5").Value = "  +0.97%  "

# Row 6
$ws.Range("E6").Value = "  -0.04%  "

# Row 7
$ws.Range("D7").Value = "'0.4626"
$ws.Range("E7").Value = "  -0.09%  "

# Row 8
$ws.Range("D8").Value = "'0.3700"
$ws.Range("E8").Value = "  +0.40%  "

# Row 9
$ws.Range("D9").Value = "'0.07288"
$ws.Range("E9").Value = "  -0.54%  "

# Row 10
$ws.Range("D10").Value = "'0.8872"
$ws.Range("E10").Value = "  +1.38%  "

# Row 11
$ws.Range("D11").Value = "'20.00"
$ws.Range("E11").Value = "  +1.89%  "

# Row 12
$ws.Range("D12").Value = "'0.07872"
$ws.Range("E12").Value = "  -0.06%  "

# Row 13
$ws.Range("D13").Value = "1.817.05"
$ws.Range("E13").Value = "  +1.32%  "

# Row 14
$ws.Range("D14").Value = "'5.392"
$ws.Range("E14").Value = "  +1.15%  "

# Row 15
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "'91.82"
$ws.Range("E15").Value = "  +0.70%  "

# Row 16
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'6.499"
$ws.Range("E16").Value = "  -0.50%  "

# Row 17
$ws.Range("E17").Value = "  -0.08%  "

# Row 18
$ws.Range("D18").Value = "'0.000008897"
$ws.Range("E18").Value = "  +0.15%  "

# Row 19
$ws.Range("E19").Value = "  +0.12%  "

# Row 20
$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").Value = "27.243.62"
$ws.Range("E20").Value = "  +0.95%  "

# Row 21
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'14.68"
$ws.Range("E21").Value = "  -0.43%  "

# Row 22
$ws.Range("D22").Value = "'5.062"
$ws.Range("E22").Value = "  -0.61%  "

# Row 23
$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").Value = "2.063.97"
$ws.Range("E23").Value = "  +3.55%  "

# Row 24
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "'10.51"
$ws.Range("E24").Value = "  -0.06%  "

# Row 25
$ws.Range("D25").Value = "'2.029"
$ws.Range("E25").Value = "  +10.13%  "

# Row 26
$ws.Range("D26").Value = "'152.00"
$ws.Range("E26").Value = "  -0.22%  "

# Row 27
$ws.Range("D27").Value = "'18.37"
$ws.Range("E27").Value = "  -0.03%  "

# Row 28
$ws.Range("D28").Value = "'2.028"
$ws.Range("E28").Value = "  -0.09%  "

# Row 29
$ws.Range("D29").Value = "'115.80"
$ws.Range("E29").Value = "  +0.21%  "

# Row 30
$ws.Range("D30").Value = "'5.028"
$ws.Range("E30").Value = "  -1.10%  "

# Row 31
$ws.Range("D31").Value = "'0.08844"
$ws.Range("E31").Value = "  -0.07%  "

# Row 32
$ws.Range("D32").Value = "'3.167"
$ws.Range("E32").Value = "  +7.12%  "

# Row 33
$ws.Range("D33").Value = "'0.7675"
$ws.Range("E33").Value = "  +4.95%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.522"
$ws.Range("E34").Value = "  +2.09%  "

# Row 35
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.159"
$ws.Range("E35").Value = "  +2.80%  "

# Row 36
$ws.Range("D36").Value = "'2.728"
$ws.Range("E36").Value = "  +12.64%  "

# Row 37
$ws.Range("D37").Value = "'1.105"
$ws.Range("E37").Value = "  +3.33%  "

# Row 38
$ws.Range("D38").Value = "'0.01936"
$ws.Range("E38").Value = "  -0.11%  "

# Row 39
$ws.Range("D39").Value = "'0.05202"
$ws.Range("E39").Value = "  -0.61%  "

# Row 40
$ws.Range("D40").Value = "'2.948"
$ws.Range("E40").Value = "  +0.26%  "

# Row 41
$ws.Range("D41").Value = "'6.992"
$ws.Range("E41").Value = "  -0.03%  "

# Row 42
$ws.Range("D42").Value = "'0.5077"
$ws.Range("E42").Value = "  -0.77%  "

# Row 43
$ws.Range("D43").Value = "'0.1621"
$ws.Range("E43").Value = "  -0.06%  "

# Row 44
$ws.Range("D44").Value = "'8.448"
$ws.Range("E44").Value = "  +4.19%  "

# Row 45
$ws.Range("D45").Value = "'0.4770"
$ws.Range("E45").Value = "  -0.75%  "

# Row 46
$ws.Range("D46").Value = "'10.31"
$ws.Range("E46").Value = "  +1.38%  "

# Row 47
$ws.Range("E47").Value = "  -0.15%  "

# Row 48
$ws.Range("D48").Value = "'102.43"
$ws.Range("E48").Value = "  +0.99%  "

# Row 49
$ws.Range("D49").Value = "'1.638"
$ws.Range("E49").Value = "  +1.44%  "

# Row 50
$ws.Range("D50").Value = "'0.06197"
$ws.Range("E50").Value = "  +0.01%  "

# Row 51
$ws.Range("D51").Value = "'65.28"
$ws.Range("E51").Value = "  +1.72%  "
